$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$json = @"
{
  "template": "Your otac {OTAC}",
  "subject": "test",
  "address": "idtest@mailinator.com",
  "senderAddress": "idtest@mailinator.com"
}
"@

$ws.Range("F2").Value = $json
$ws.Range("F2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 68.65

$ws.Range("F2").Select() | Out-Null
